$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 ("Number of Promotions") - shifts rows up
$ws.Rows.Item(4).Delete()

# Row 6 ("Number of Dependents") is now row 5 after the shift above - delete it
$ws.Rows.Item(5).Delete()
